$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text: Volume number and report date range ---
$ws.Range("C8").Value = "Volume 33   Number  2"
$ws.Range("C9").Value = "Report Covering the Week  1/5/2026  Through  1/11/2026"

# --- Cells switching FROM a numeric value TO the text placeholder ---
# Donor C31 carries style 13 + shared text "0"; E31 carries style 13 + shared text "***.*"
$ws.Range("C31").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("C31").Copy()
$ws.Range("C14").PasteSpecial(-4163)
$ws.Range("C31").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C31").Copy()
$ws.Range("C20").PasteSpecial(-4163)
$ws.Range("C31").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("C31").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E31").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E31").Copy()
$ws.Range("E22").PasteSpecial(-4163)
$ws.Range("C31").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C31").Copy()
$ws.Range("C29").PasteSpecial(-4163)
$ws.Range("C31").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("C31").Copy()
$ws.Range("C30").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# --- Cells switching FROM the text placeholder TO a numeric value ---
# Determine donor by destination column's existing numeric style (14 = count, 15 = percent)
$ws.Range("F14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = 1
$ws.Range("F14").Copy()
$ws.Range("I15").PasteSpecial(-4122)
$ws.Range("I15").Value = 1
$ws.Range("F14").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").Value = 2
$ws.Range("L16").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("E16").Value = 0
$ws.Range("F14").Copy()
$ws.Range("J16").PasteSpecial(-4122)
$ws.Range("J16").Value = 2
$ws.Range("L16").Copy()
$ws.Range("K16").PasteSpecial(-4122)
$ws.Range("K16").Value = 200
$ws.Range("F14").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = 1
$ws.Range("L16").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E18").Value = 200
$ws.Range("F14").Copy()
$ws.Range("I18").PasteSpecial(-4122)
$ws.Range("I18").Value = 3
$ws.Range("F14").Copy()
$ws.Range("J18").PasteSpecial(-4122)
$ws.Range("J18").Value = 1
$ws.Range("L16").Copy()
$ws.Range("K18").PasteSpecial(-4122)
$ws.Range("K18").Value = 200
$ws.Range("F14").Copy()
$ws.Range("I19").PasteSpecial(-4122)
$ws.Range("I19").Value = 2
$ws.Range("F14").Copy()
$ws.Range("J20").PasteSpecial(-4122)
$ws.Range("J20").Value = 1
$ws.Range("L16").Copy()
$ws.Range("K20").PasteSpecial(-4122)
$ws.Range("K20").Value = -100
$ws.Range("L16").Copy()
$ws.Range("M20").PasteSpecial(-4122)
$ws.Range("M20").Value = -100
$ws.Range("F14").Copy()
$ws.Range("J23").PasteSpecial(-4122)
$ws.Range("J23").Value = 1
$ws.Range("L16").Copy()
$ws.Range("K23").PasteSpecial(-4122)
$ws.Range("K23").Value = 200
$ws.Range("L16").Copy()
$ws.Range("M24").PasteSpecial(-4122)
$ws.Range("M24").Value = 283.333333333333
$ws.Range("L16").Copy()
$ws.Range("L25").PasteSpecial(-4122)
$ws.Range("L25").Value = 150
$ws.Range("F14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = 1
$ws.Range("F14").Copy()
$ws.Range("I27").PasteSpecial(-4122)
$ws.Range("I27").Value = 1
$ws.Range("F14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C28").Value = 2
$ws.Range("F14").Copy()
$ws.Range("I28").PasteSpecial(-4122)
$ws.Range("I28").Value = 2
$ws.Range("L16").Copy()
$ws.Range("L29").PasteSpecial(-4122)
$ws.Range("L29").Value = -100
$ws.Range("L16").Copy()
$ws.Range("N29").PasteSpecial(-4122)
$ws.Range("N29").Value = -100
$ws.Range("L16").Copy()
$ws.Range("L30").PasteSpecial(-4122)
$ws.Range("L30").Value = -100
$ws.Range("L16").Copy()
$ws.Range("N30").PasteSpecial(-4122)
$ws.Range("N30").Value = -100
$excel.CutCopyMode = $false

# --- Plain numeric updates (style unchanged) ---
$ws.Range("F15").Value = 3
$ws.Range("H15").Value = 0
$ws.Range("L15").Value = -50
$ws.Range("C16").Value = 2
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = 83.333333333333
$ws.Range("I16").Value = 6
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = -60
$ws.Range("C17").Value = 12
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 39
$ws.Range("G17").Value = 28
$ws.Range("H17").Value = 39.285714285714
$ws.Range("I17").Value = 17
$ws.Range("J17").Value = 8
$ws.Range("K17").Value = 112.5
$ws.Range("L17").Value = 30.76923076923
$ws.Range("M17").Value = 240
$ws.Range("N17").Value = 41.666666666666
$ws.Range("C18").Value = 3
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = 16.666666666666
$ws.Range("L18").Value = 50
$ws.Range("M18").Value = -57.142857142857
$ws.Range("N18").Value = -80
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 16
$ws.Range("G19").Value = 21
$ws.Range("H19").Value = -23.809523809523
$ws.Range("J19").Value = 8
$ws.Range("K19").Value = -75
$ws.Range("L19").Value = -60
$ws.Range("M19").Value = -33.333333333333
$ws.Range("N19").Value = -83.333333333333
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = -50
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 14
$ws.Range("E21").Value = 42.857142857142
$ws.Range("F21").Value = 80
$ws.Range("G21").Value = 70
$ws.Range("H21").Value = 14.285714285714
$ws.Range("I21").Value = 29
$ws.Range("J21").Value = 20
$ws.Range("K21").Value = 45
$ws.Range("L21").Value = 3.571428571428
$ws.Range("M21").Value = 20.833333333333
$ws.Range("N21").Value = -60.273972602739
$ws.Range("C23").Value = 1
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 10
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = 100
$ws.Range("I23").Value = 3
$ws.Range("L23").Value = -25
$ws.Range("M23").Value = 50
$ws.Range("C24").Value = 14
$ws.Range("D24").Value = 5
$ws.Range("E24").Value = 180
$ws.Range("F24").Value = 71
$ws.Range("G24").Value = 35
$ws.Range("H24").Value = 102.857142857143
$ws.Range("I24").Value = 23
$ws.Range("J24").Value = 12
$ws.Range("K24").Value = 91.666666666666
$ws.Range("L24").Value = 43.75
$ws.Range("C25").Value = 2
$ws.Range("F25").Value = 16
$ws.Range("H25").Value = 1500
$ws.Range("I25").Value = 5
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 9
$ws.Range("E26").Value = -22.222222222222
$ws.Range("F26").Value = 39
$ws.Range("H26").Value = -4.878048780487
$ws.Range("I26").Value = 10
$ws.Range("J26").Value = 13
$ws.Range("K26").Value = -23.076923076923
$ws.Range("L26").Value = -52.380952380952
$ws.Range("M26").Value = -23.076923076923
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = 0
$ws.Range("L27").Value = -50
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 100
$ws.Range("J28").Value = 2
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 100
